$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) First paragraph: split "This is a Microsoft word document." into four
#    runs - the original text (now with two trailing spaces) plus a new
#    red (C00000) parenthetical comment appended across three runs.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)

$para1Xml = '<w:p ' + $wns + '>' +
  '<w:r><w:t xml:space="preserve">This is a Microsoft word document.  </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>(This is a change ' + [char]0x2013 + ' Ve</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>rsion for branch alternate</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>)</w:t></w:r>' +
  '</w:p>'

$p1.Range.InsertXML($para1Xml)

# ---------------------------------------------------------------------------
# 2) Insert a brand-new, empty paragraph right after
#    "It will be treated as a binary file by Git." carrying shading +
#    paragraph-mark run formatting but no text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "It will be treated as a binary file by Git.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "It will be treated as a binary file by Git.^p", 2) | Out-Null

$newPara = $d.Paragraphs.Item(3)

$newParaXml = '<w:p ' + $wns + '>' +
  '<w:pPr>' +
    '<w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/>' +
    '<w:rPr>' +
      '<w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>' +
      '<w:b/>' +
      '<w:bCs/>' +
      '<w:color w:val="202122"/>' +
    '</w:rPr>' +
  '</w:pPr>' +
  '</w:p>'

$newPara.Range.InsertXML($newParaXml)

Write-Host "edit complete"
